$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 295.1111
$ws.Range("I9").Value = 259.33334
$ws.Range("J9").Value = 366.66666
$ws.Range("K9").Value = 259.33334
$ws.Range("L9").Value = 366.66666
$ws.Range("M9").Value = -90.33334000000002
$ws.Range("N9").Value = -704.66666

$ws.Range("H19").Value = 2071.4285
$ws.Range("I19").Value = 1750
$ws.Range("J19").Value = 2500
$ws.Range("K19").Value = 1750
$ws.Range("L19").Value = 2500
$ws.Range("M19").Value = -1575
$ws.Range("N19").Value = -2850

$ws.Range("H116").Value = 3666.6667
$ws.Range("I116").Value = 3502.5
$ws.Range("J116").Value = 3995
$ws.Range("K116").Value = 3502.5
$ws.Range("L116").Value = 3995
$ws.Range("M116").Value = -60.5
$ws.Range("N116").Value = -10879

$ws.Range("N130").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("N130").Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2223
$ws.Range("I61").Value = 2155.64
$ws.Range("J61").Value = 2503.6667
$ws.Range("K61").Value = 2155.64
$ws.Range("L61").Value = 2503.6667
$ws.Range("M61").Value = -1943.64
$ws.Range("N61").Value = -2927.6667

$ws.Range("H64").Value = 25000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 25000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25496

$ws.Range("H67").Value = 25000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 25000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26716

$ws.Range("H74").Value = 3513.8333
$ws.Range("I74").Value = 3900.7778
$ws.Range("J74").Value = 2353
$ws.Range("K74").Value = 3900.7778
$ws.Range("L74").Value = 2353
$ws.Range("M74").Value = -3026.7778
$ws.Range("N74").Value = -4101

$ws.Range("H77").Value = 3513.8333
$ws.Range("I77").Value = 3900.7778
$ws.Range("J77").Value = 2353
$ws.Range("K77").Value = 19503.889
$ws.Range("L77").Value = 11765
$ws.Range("M77").Value = -15135.889
$ws.Range("N77").Value = -20501

$ws.Range("H97").Value = 1357.1875
$ws.Range("I97").Value = 747.3
$ws.Range("J97").Value = 10505.5
$ws.Range("K97").Value = 747.3
$ws.Range("L97").Value = 10505.5
$ws.Range("M97").Value = -251.3
$ws.Range("N97").Value = -11497.5

$ws.Range("H102").Value = 1943.4
$ws.Range("I102").Value = 1991.4445
$ws.Range("J102").Value = 1511
$ws.Range("K102").Value = 1991.4445
$ws.Range("L102").Value = 1511
$ws.Range("M102").Value = -369.4445000000001
$ws.Range("N102").Value = -4755

$ws.Range("H136").Value = 2223
$ws.Range("I136").Value = 2155.64
$ws.Range("J136").Value = 2503.6667
$ws.Range("K136").Value = 6466.92
$ws.Range("L136").Value = 7511.000100000001
$ws.Range("M136").Value = -3916.92
$ws.Range("N136").Value = -12611.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 45000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 45000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46372

$ws.Range("H65").Value = 45000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 45000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141864

$ws.Range("H86").Value = 2661.4167
$ws.Range("I86").Value = 1944.7142
$ws.Range("J86").Value = 4591
$ws.Range("K86").Value = 1944.7142
$ws.Range("L86").Value = 4591
$ws.Range("M86").Value = -821.7141999999999
$ws.Range("N86").Value = -6837

$ws.Range("H89").Value = 2661.4167
$ws.Range("I89").Value = 1944.7142
$ws.Range("J89").Value = 4591
$ws.Range("K89").Value = 9723.571
$ws.Range("L89").Value = 22955
$ws.Range("M89").Value = -4107.571
$ws.Range("N89").Value = -34187

$ws.Range("H94").Value = 2328.2917
$ws.Range("I94").Value = 2408.8696
$ws.Range("J94").Value = 475
$ws.Range("K94").Value = 2408.8696
$ws.Range("L94").Value = 475
$ws.Range("M94").Value = -1957.8696
$ws.Range("N94").Value = -1377

$ws.Range("H99").Value = 46078.332
$ws.Range("I99").Value = 51525.625
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 51525.625
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -50027.625
$ws.Range("N99").Value = -5496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2243.5862
$ws.Range("I31").Value = 2128.5715
$ws.Range("J31").Value = 2545.5
$ws.Range("K31").Value = 2128.5715
$ws.Range("L31").Value = 2545.5
$ws.Range("M31").Value = -1833.5715
$ws.Range("N31").Value = -3135.5

$ws.Range("H34").Value = 2243.5862
$ws.Range("I34").Value = 2128.5715
$ws.Range("J34").Value = 2545.5
$ws.Range("K34").Value = 2128.5715
$ws.Range("L34").Value = 2545.5
$ws.Range("M34").Value = -1926.5715
$ws.Range("N34").Value = -2949.5

$ws.Range("H99").Value = 32062138
$ws.Range("I99").Value = 12196378
$ws.Range("J99").Value = 40008444
$ws.Range("K99").Value = 12196378
$ws.Range("L99").Value = 40008444
$ws.Range("M99").Value = -12194880
$ws.Range("N99").Value = -40011440

$ws.Range("H107").Value = 10410.286
$ws.Range("I107").Value = 652.2
$ws.Range("J107").Value = 19281.273
$ws.Range("K107").Value = 652.2
$ws.Range("L107").Value = 19281.273
$ws.Range("M107").Value = 1267.8
$ws.Range("N107").Value = -23121.273

$ws.Range("H126").Value = 32062138
$ws.Range("I126").Value = 12196378
$ws.Range("J126").Value = 40008444
$ws.Range("K126").Value = 36589134
$ws.Range("L126").Value = 120025332
$ws.Range("M126").Value = -36586664
$ws.Range("N126").Value = -120030272

$ws.Range("H134").Value = 2127.5454
$ws.Range("I134").Value = 1726.4117
$ws.Range("J134").Value = 3491.4
$ws.Range("K134").Value = 5179.2351
$ws.Range("L134").Value = 10474.2
$ws.Range("M134").Value = -2644.2351
$ws.Range("N134").Value = -15544.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 787.4
$ws.Range("I86").Value = 253
$ws.Range("J86").Value = 1589
$ws.Range("K86").Value = 759
$ws.Range("L86").Value = 4767
$ws.Range("M86").Value = 427
$ws.Range("N86").Value = -7139

$ws.Range("H89").Value = 787.4
$ws.Range("I89").Value = 253
$ws.Range("J89").Value = 1589
$ws.Range("K89").Value = 2277
$ws.Range("L89").Value = 14301
$ws.Range("M89").Value = 3651
$ws.Range("N89").Value = -26157

$ws.Range("H98").Value = 348.5
$ws.Range("I98").Value = 348.5
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1045.5
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 452.5

$ws.Range("H122").Value = 854.46155
$ws.Range("I122").Value = 944
$ws.Range("J122").Value = 750
$ws.Range("K122").Value = 8496
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -6046
$ws.Range("N122").Value = -11650

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 17496.166
$ws.Range("I97").Value = 994.3333
$ws.Range("J97").Value = 33998
$ws.Range("K97").Value = 994.3333
$ws.Range("L97").Value = 33998
$ws.Range("M97").Value = -498.3333
$ws.Range("N97").Value = -34990

$ws.Range("H126").Value = 2710.8125
$ws.Range("I126").Value = 2762.3572
$ws.Range("J126").Value = 2350
$ws.Range("K126").Value = 8287.071599999999
$ws.Range("L126").Value = 7050
$ws.Range("M126").Value = -5817.071599999999
$ws.Range("N126").Value = -11990

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10981.818
$ws.Range("I7").Value = 13591.667
$ws.Range("J7").Value = 7850
$ws.Range("K7").Value = 13591.667
$ws.Range("L7").Value = 7850
$ws.Range("M7").Value = -13479.667
$ws.Range("N7").Value = -8074

$ws.Range("H40").Value = 23350.8
$ws.Range("I40").Value = 27626
$ws.Range("J40").Value = 6250
$ws.Range("K40").Value = 27626
$ws.Range("L40").Value = 6250
$ws.Range("M40").Value = -27490
$ws.Range("N40").Value = -6522

$ws.Range("H61").Value = 3387.7778
$ws.Range("I61").Value = 2902.6
$ws.Range("J61").Value = 3994.25
$ws.Range("K61").Value = 2902.6
$ws.Range("L61").Value = 3994.25
$ws.Range("M61").Value = -2700.6
$ws.Range("N61").Value = -4398.25

$ws.Range("H93").Value = 2936.3635
$ws.Range("I93").Value = 3462.5
$ws.Range("J93").Value = 1533.3334
$ws.Range("K93").Value = 3462.5
$ws.Range("L93").Value = 1533.3334
$ws.Range("M93").Value = -2214.5
$ws.Range("N93").Value = -4029.3334

$ws.Range("H100").Value = 61464.65
$ws.Range("I100").Value = 114760.8
$ws.Range("J100").Value = 8168.5
$ws.Range("K100").Value = 114760.8
$ws.Range("L100").Value = 8168.5
$ws.Range("M100").Value = -114219.8
$ws.Range("N100").Value = -9250.5

$ws.Range("H113").Value = 3387.7778
$ws.Range("I113").Value = 2902.6
$ws.Range("J113").Value = 3994.25
$ws.Range("K113").Value = 2902.6
$ws.Range("L113").Value = 3994.25
$ws.Range("M113").Value = -732.5999999999999
$ws.Range("N113").Value = -8334.25

$ws.Range("H126").Value = 10981.818
$ws.Range("I126").Value = 13591.667
$ws.Range("J126").Value = 7850
$ws.Range("K126").Value = 40775.001
$ws.Range("L126").Value = 23550
$ws.Range("M126").Value = -38305.001
$ws.Range("N126").Value = -28490

$ws.Range("H132").Value = 3665.625
$ws.Range("I132").Value = 2235.8064
$ws.Range("J132").Value = 8590.556
$ws.Range("K132").Value = 6707.4192
$ws.Range("L132").Value = 25771.668
$ws.Range("M132").Value = -4177.4192
$ws.Range("N132").Value = -30831.668

$ws.Range("H136").Value = 2927.6365
$ws.Range("I136").Value = 2507.9
$ws.Range("J136").Value = 7125
$ws.Range("K136").Value = 7523.700000000001
$ws.Range("L136").Value = 21375
$ws.Range("M136").Value = -4973.700000000001
$ws.Range("N136").Value = -26475

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1866.5405
$ws.Range("I136").Value = 1114.32
$ws.Range("J136").Value = 3433.6667
$ws.Range("K136").Value = 3342.96
$ws.Range("L136").Value = 10301.0001
$ws.Range("M136").Value = -792.96
$ws.Range("N136").Value = -15401.0001
